# Phase 5 program selection flow - data migration script
# Adds new course rows, track metadata columns, bucket definitions,
# and course-to-bucket mappings described in the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: map column letters (A, B, C, ...) used in the data tables below to
# 1-based column indices, and write a cell only when the hashtable actually
# supplies a value for that column (so we don't clobber anything with blanks).
# ---------------------------------------------------------------------------
function Get-ColIndex($letter) {
    $idx = 0
    foreach ($ch in $letter.ToCharArray()) {
        $idx = $idx * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $idx
}

function Set-RowValues($ws, $rowData, $cols) {
    $r = $rowData.Row
    foreach ($col in $cols) {
        if ($rowData.ContainsKey($col)) {
            $c = Get-ColIndex $col
            $ws.Cells.Item($r, $c).Value = $rowData[$col]
        }
    }
}

# ---------------------------------------------------------------------------
# 1) courses sheet - append rows 46-65 (new course catalog entries)
# ---------------------------------------------------------------------------
$courseRows = @(
    @{ Row = 46; A = 'REAL 3001'; B = 'Introduction to Commercial Real Estate'; C = 3; D = 3000; E = 'ECON 1103; ACCO 1030'; F = 'ECON 1103;ACCO 1030'; H = $true; I = $true; J = $true; K = '2025 Fall; 2025 Summer; 2025 Spring; 2024 Fall'; L = 'high' }
    @{ Row = 47; A = 'REAL 4061'; B = 'Real Estate Modeling'; C = 3; D = 4000; E = 'REAL 3001 or REAL 4002 or FINA 4002'; F = 'REAL 3001 or FINA 4002'; H = $false; I = $true; J = $false; K = '2025 Spring; 2024 Spring; 2023 Spring; 2022 Spring'; L = 'high' }
    @{ Row = 48; A = 'ACCO 4040'; B = 'International Accounting'; C = 3; D = 4000; E = 'ACCO 4020 (may be concurrent)'; F = 'ACCO 4020'; G = 'may_be_concurrent'; H = $true; I = $false; J = $false; K = '2025 Fall; 2024 Fall; 2023 Fall; 2023 Spring'; L = 'medium' }
    @{ Row = 49; A = 'ECON 4040'; B = 'International Economics'; C = 3; D = 4000; E = 'ECON 1103; ECON 1104'; F = 'ECON 1103;ECON 1104'; H = $true; I = $true; J = $false; K = '2025 Fall; 2024 Spring; 2023 Spring; 2021 Fall'; L = 'medium' }
    @{ Row = 50; A = 'ECON 4044'; B = 'Global Integration of Financial Sectors'; C = 3; D = 4000; E = 'ECON 1103; ECON 1104'; F = 'ECON 1103;ECON 1104'; H = $false; I = $false; J = $false; K = '2018 Fall; 2018 Spring; 2017 Spring; 2016 Fall'; L = 'low'; M = 'Not offered since 2018' }
    @{ Row = 51; A = 'FINAI 4931'; B = 'Topics in Finance-International'; C = 3; D = 4000; E = 'FINA 3001; OIE consent'; F = 'FINA 3001'; G = 'enrollment_requirement'; H = $false; I = $true; J = $false; K = '2025 Spring; 2024 Spring'; L = 'medium'; M = 'Requires OIE consent' }
    @{ Row = 52; A = 'INBUI 4931'; B = 'Topics in International Business-International'; C = 3; D = 4000; E = 'OIE consent'; F = 'none'; G = 'enrollment_requirement'; H = $false; I = $true; J = $true; K = '2025 Summer; 2025 Spring; 2024 Summer; 2024 Spring'; L = 'high'; M = 'Requires OIE consent' }
    @{ Row = 53; A = 'MARK 4040'; B = 'International Marketing'; C = 3; D = 4000; E = 'MARK 3001'; F = 'MARK 3001'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high' }
    @{ Row = 54; A = 'ACCO 3001'; B = 'Intermediate Accounting I'; C = 3; D = 3000; E = 'ACCO 1031'; F = 'ACCO 1031'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high' }
    @{ Row = 55; A = 'ACCO 4020'; B = 'Intermediate Accounting II'; C = 3; D = 4000; E = 'ACCO 3001'; F = 'ACCO 3001'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high' }
    @{ Row = 56; A = 'ACCO 4080'; B = 'Analysis of Financial Statements'; C = 3; D = 4000; E = 'ACCO 3001 or AIM/CBP admission'; F = 'ACCO 3001'; G = 'admitted_program'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high'; M = 'Also accepts AIM/CBP admission' }
    @{ Row = 57; A = 'ECON 4080'; B = 'Money, Banking and Monetary Policy'; C = 3; D = 4000; E = 'ECON 3004'; F = 'ECON 3004'; H = $true; I = $true; J = $true; K = '2024 Fall; 2023 Summer; 2023 Spring; 2022 Spring'; L = 'low'; M = 'Irregular offering schedule' }
    @{ Row = 58; A = 'INSY 4051'; B = 'Business Applications Development'; C = 3; D = 4000; E = 'INSY 3001 or ACCO 4050'; F = 'INSY 3001 or ACCO 4050'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high' }
    @{ Row = 59; A = 'INSY 4053'; B = 'Project Management'; C = 3; D = 4000; E = 'INSY 3001 or ACCO 4050 or instructor consent'; F = 'INSY 3001 or ACCO 4050'; G = 'instructor_consent'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high' }
    @{ Row = 60; A = 'MARK 4094'; B = 'Professional Selling'; C = 3; D = 4000; E = 'MARK 3001'; F = 'MARK 3001'; H = $true; I = $true; J = $true; K = '2025 Fall; 2025 Summer; 2025 Spring; 2024 Fall'; L = 'high' }
    @{ Row = 61; A = 'REAL 4100'; B = 'Commercial Real Estate Development'; C = 3; D = 4000; E = 'REAL 3001 or REAL 4002 or FINA 4002'; F = 'REAL 3001 or FINA 4002'; H = $false; I = $true; J = $false; K = '2025 Spring; 2024 Spring; 2023 Spring; 2022 Spring'; L = 'high' }
    @{ Row = 62; A = 'MARK 3001'; B = 'Introduction to Marketing'; C = 3; D = 3000; E = 'Soph. stndg.; ECON 1001 or ECON 1103'; F = 'ECON 1103'; G = 'standing_requirement'; H = $true; I = $true; J = $true; K = '2025 Fall; 2025 Summer; 2025 Spring; 2024 Fall'; L = 'high' }
    @{ Row = 63; A = 'ECON 3004'; B = 'Intermediate Macroeconomic Analysis'; C = 3; D = 3000; E = 'ECON 1103; ECON 1104; MATH 1400 or equiv.'; F = 'ECON 1103;ECON 1104;MATH 1400'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Spring; 2024 Fall; 2024 Spring'; L = 'high' }
    @{ Row = 64; A = 'INSY 3001'; B = 'Introduction to Information Systems'; C = 3; D = 3000; E = 'Soph. stndg.'; F = 'none'; G = 'standing_requirement'; H = $true; I = $true; J = $true; K = '2025 Fall; 2025 Summer; 2025 Spring; 2024 Fall'; L = 'high' }
    @{ Row = 65; A = 'ACCO 4050'; B = 'Accounting Information Systems'; C = 3; D = 4000; E = 'ACCO 1031'; F = 'ACCO 1031'; H = $true; I = $true; J = $false; K = '2025 Fall; 2025 Summer; 2024 Fall; 2024 Summer'; L = 'high' }
)

$coursesWs = $wb.Worksheets.Item("courses")
$courseCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")
foreach ($row in $courseRows) {
    Set-RowValues $coursesWs $row $courseCols
}

# ---------------------------------------------------------------------------
# 2) tracks sheet - add `kind` / `parent_major_id` columns (D, E)
# ---------------------------------------------------------------------------
$tracksWs = $wb.Worksheets.Item("tracks")

$tracksWs.Cells.Item(1,4).Value = "kind"
$tracksWs.Cells.Item(1,5).Value = "parent_major_id"

$tracksWs.Cells.Item(2,4).Value = "major"
$tracksWs.Cells.Item(2,5).Value = ""

$tracksWs.Cells.Item(3,4).Value = "track"
$tracksWs.Cells.Item(3,5).Value = "FIN_MAJOR"

$tracksWs.Cells.Item(4,4).Value = "track"
$tracksWs.Cells.Item(4,5).Value = "FIN_MAJOR"

# ---------------------------------------------------------------------------
# 3) buckets sheet - append rows 7-10 (new concentration bucket definitions)
# ---------------------------------------------------------------------------
$bucketRows = @(
    @{ Row = 7; A = 'FP_CONC'; B = 'FP_CORE'; C = 'Financial Planning Required'; D = 1; E = 7; H = $false; I = 'core' }
    @{ Row = 8; A = 'CB_CONC'; B = 'CB_CORE'; C = 'Commercial Banking Required'; D = 1; E = 9; H = $false; I = 'core' }
    @{ Row = 9; A = 'CB_CONC'; B = 'CB_INTL'; C = 'International Requirement'; D = 2; E = 1; H = $true }
    @{ Row = 10; A = 'CB_CONC'; B = 'CB_ELEC'; C = 'Commercial Banking Elective'; D = 3; E = 2; H = $true; I = 'elective' }
)

$bucketsWs = $wb.Worksheets.Item("buckets")
$bucketCols = @("A","B","C","D","E","F","G","H","I")
foreach ($row in $bucketRows) {
    Set-RowValues $bucketsWs $row $bucketCols
}

# ---------------------------------------------------------------------------
# 4) course_bucket sheet - append rows 98-136 (course -> bucket mappings for
#    the new Financial Planning / Commercial Banking concentrations)
# ---------------------------------------------------------------------------
$courseBucketRows = @(
    @{ Row = 98; A = 'FP_CONC'; B = 'FINA 3001'; C = 'FP_CORE' }
    @{ Row = 99; A = 'FP_CONC'; B = 'FINA 4001'; C = 'FP_CORE' }
    @{ Row = 100; A = 'FP_CONC'; B = 'FINA 4011'; C = 'FP_CORE' }
    @{ Row = 101; A = 'FP_CONC'; B = 'FINA 4020'; C = 'FP_CORE' }
    @{ Row = 102; A = 'FP_CONC'; B = 'FINA 4121'; C = 'FP_CORE' }
    @{ Row = 103; A = 'FP_CONC'; B = 'FINA 4122'; C = 'FP_CORE' }
    @{ Row = 104; A = 'FP_CONC'; B = 'FINA 4123'; C = 'FP_CORE' }
    @{ Row = 105; A = 'CB_CONC'; B = 'FINA 3001'; C = 'CB_CORE' }
    @{ Row = 106; A = 'CB_CONC'; B = 'FINA 3002'; C = 'CB_CORE' }
    @{ Row = 107; A = 'CB_CONC'; B = 'FINA 4001'; C = 'CB_CORE' }
    @{ Row = 108; A = 'CB_CONC'; B = 'FINA 4011'; C = 'CB_CORE' }
    @{ Row = 109; A = 'CB_CONC'; B = 'FINA 4050'; C = 'CB_CORE' }
    @{ Row = 110; A = 'CB_CONC'; B = 'REAL 4061'; C = 'CB_CORE' }
    @{ Row = 111; A = 'CB_CONC'; B = 'FINA 4210'; C = 'CB_CORE' }
    @{ Row = 112; A = 'CB_CONC'; B = 'FINA 4211'; C = 'CB_CORE' }
    @{ Row = 113; A = 'CB_CONC'; B = 'FINA 4212'; C = 'CB_CORE' }
    @{ Row = 114; A = 'CB_CONC'; B = 'REAL 3001'; C = 'CB_CORE' }
    @{ Row = 115; A = 'CB_CONC'; B = 'ACCO 4040'; C = 'CB_INTL' }
    @{ Row = 116; A = 'CB_CONC'; B = 'ECON 4040'; C = 'CB_INTL' }
    @{ Row = 117; A = 'CB_CONC'; B = 'ECON 4044'; C = 'CB_INTL' }
    @{ Row = 118; A = 'CB_CONC'; B = 'FINA 4040'; C = 'CB_INTL' }
    @{ Row = 119; A = 'CB_CONC'; B = 'FINAI 4931'; C = 'CB_INTL' }
    @{ Row = 120; A = 'CB_CONC'; B = 'INBUI 4931'; C = 'CB_INTL' }
    @{ Row = 121; A = 'CB_CONC'; B = 'MARK 4040'; C = 'CB_INTL' }
    @{ Row = 122; A = 'CB_CONC'; B = 'ACCO 3001'; C = 'CB_ELEC' }
    @{ Row = 123; A = 'CB_CONC'; B = 'ACCO 4020'; C = 'CB_ELEC' }
    @{ Row = 124; A = 'CB_CONC'; B = 'ACCO 4080'; C = 'CB_ELEC' }
    @{ Row = 125; A = 'CB_CONC'; B = 'FINA 4002'; C = 'CB_ELEC' }
    @{ Row = 126; A = 'CB_CONC'; B = 'FINA 4065'; C = 'CB_ELEC' }
    @{ Row = 127; A = 'CB_CONC'; B = 'FINA 4075'; C = 'CB_ELEC' }
    @{ Row = 128; A = 'CB_CONC'; B = 'FINA 4081'; C = 'CB_ELEC' }
    @{ Row = 129; A = 'CB_CONC'; B = 'FINA 4082'; C = 'CB_ELEC' }
    @{ Row = 130; A = 'CB_CONC'; B = 'FINA 4084'; C = 'CB_ELEC' }
    @{ Row = 131; A = 'CB_CONC'; B = 'ECON 4080'; C = 'CB_ELEC' }
    @{ Row = 132; A = 'CB_CONC'; B = 'INSY 4051'; C = 'CB_ELEC' }
    @{ Row = 133; A = 'CB_CONC'; B = 'INSY 4053'; C = 'CB_ELEC' }
    @{ Row = 134; A = 'CB_CONC'; B = 'MARK 4094'; C = 'CB_ELEC' }
    @{ Row = 135; A = 'CB_CONC'; B = 'REAL 4100'; C = 'CB_ELEC' }
    @{ Row = 136; A = 'CB_CONC'; B = 'REAL 4061'; C = 'CB_ELEC' }
)

$courseBucketWs = $wb.Worksheets.Item("course_bucket")
$courseBucketCols = @("A","B","C")
foreach ($row in $courseBucketRows) {
    Set-RowValues $courseBucketWs $row $courseBucketCols
}
